$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.264496
$ws.Range("H2").Value = 0.793488
$ws.Range("I2").Value = 0.001006353962629067
$ws.Range("J2").Value = 0.001006353962629067
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.1148133333333333
$ws.Range("N2").Value = 0.34444
$ws.Range("O2").Value = 0.03343792635928704
$ws.Range("P2").Value = 0.03343792635928704
$ws.Range("Q2").Value = 0.03036766741333333
$ws.Range("R2").Value = 0.27330900672
$ws.Range("S2").Value = 0.00003365038969376745
$ws.Range("T2").Value = 0.00003365038969376745
$ws.Range("G3").Value = 0.264496
$ws.Range("H3").Value = 0.793488
$ws.Range("I3").Value = 0.001006353962629067
$ws.Range("J3").Value = 0.001006353962629067
$ws.Range("O3").Value = 0.9249645515654102
$ws.Range("P3").Value = 0.9249645515654102
$ws.Range("Q3").Value = 0.8400346232373335
$ws.Range("R3").Value = 7.560311609136001
$ws.Range("S3").Value = 0.0009308417417592687
$ws.Range("T3").Value = 0.0009308417417592686
$ws.Range("G4").Value = 0.264496
$ws.Range("H4").Value = 0.793488
$ws.Range("I4").Value = 0.001006353962629067
$ws.Range("J4").Value = 0.001006353962629067
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.1428303333333333
$ws.Range("N4").Value = 0.428491
$ws.Range("O4").Value = 0.04159752207530271
$ws.Range("P4").Value = 0.04159752207530271
$ws.Range("Q4").Value = 0.03777805184533334
$ws.Range("R4").Value = 0.340002466608
$ws.Range("S4").Value = 0.00004186183117603097
$ws.Range("T4").Value = 0.00004186183117603097
$ws.Range("H5").Value = 737.537796
$ws.Range("I5").Value = 0.9353942133886188
$ws.Range("J5").Value = 0.9353942133886189
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.1148133333333333
$ws.Range("N5").Value = 0.34444
$ws.Range("O5").Value = 0.03343792635928704
$ws.Range("P5").Value = 0.03343792635928704
$ws.Range("Q5").Value = 28.22639093936
$ws.Range("R5").Value = 254.03751845424
$ws.Range("S5").Value = 0.03127764282419187
$ws.Range("T5").Value = 0.03127764282419187
$ws.Range("H6").Value = 737.537796
$ws.Range("I6").Value = 0.9353942133886188
$ws.Range("J6").Value = 0.9353942133886189
$ws.Range("Q6").Value = 780.802336753868
$ws.Range("R6").Value = 7027.221030784813
$ws.Range("S6").Value = 0.8652064891238834
$ws.Range("T6").Value = 0.8652064891238835
$ws.Range("H7").Value = 737.537796
$ws.Range("I7").Value = 0.9353942133886188
$ws.Range("J7").Value = 0.9353942133886189
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.1428303333333333
$ws.Range("N7").Value = 0.428491
$ws.Range("O7").Value = 0.04159752207530271
$ws.Range("P7").Value = 0.04159752207530271
$ws.Range("Q7").Value = 35.114256416204
$ws.Range("R7").Value = 316.028307745836
$ws.Range("S7").Value = 0.03891008144054348
$ws.Range("T7").Value = 0.03891008144054349
$ws.Range("G8").Value = 16.71558533333333
$ws.Range("H8").Value = 50.146756
$ws.Range("I8").Value = 0.06359943264875202
$ws.Range("J8").Value = 0.06359943264875202
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.1148133333333333
$ws.Range("N8").Value = 0.34444
$ws.Range("O8").Value = 0.03343792635928704
$ws.Range("P8").Value = 0.03343792635928704
$ws.Range("Q8").Value = 1.919172070737778
$ws.Range("R8").Value = 17.27254863664
$ws.Range("S8").Value = 0.002126633145401406
$ws.Range("T8").Value = 0.002126633145401406
$ws.Range("G9").Value = 16.71558533333333
$ws.Range("H9").Value = 50.146756
$ws.Range("I9").Value = 0.06359943264875202
$ws.Range("J9").Value = 0.06359943264875202
$ws.Range("O9").Value = 0.9249645515654102
$ws.Range("P9").Value = 0.9249645515654102
$ws.Range("Q9").Value = 53.08840370999245
$ws.Range("R9").Value = 477.795633389932
$ws.Range("S9").Value = 0.05882722069976742
$ws.Range("T9").Value = 0.05882722069976742
$ws.Range("G10").Value = 16.71558533333333
$ws.Range("H10").Value = 50.146756
$ws.Range("I10").Value = 0.06359943264875202
$ws.Range("J10").Value = 0.06359943264875202
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.1428303333333333
$ws.Range("N10").Value = 0.428491
$ws.Range("O10").Value = 0.04159752207530271
$ws.Range("P10").Value = 0.04159752207530271
$ws.Range("Q10").Value = 2.387492625021778
$ws.Range("R10").Value = 21.487433625196
$ws.Range("S10").Value = 0.00264557880358319
$ws.Range("T10").Value = 0.00264557880358319
